$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 'it''s a bad experience for him'
$ws.Range("E2").Value = 'and now he''s moved out I don''t have my job he''s moved out you know he''s got the kids he said I can''t be trusted with my kids so you don''t'
$ws.Range("F2").Value = 'oh okay so if you were to move out would you still be able to work?'

$ws.Range("D3").Value = 'something else could change'
$ws.Range("E3").Value = 'yeah I don''t know what'
$ws.Range("F3").Value = 'do you want to talk about it?'

$ws.Range("D4").Value = 'and why why do you want to lose weight'
$ws.Range("E4").Value = 'well so really recently actually my older sister has had a lot of health concerns relate to her weight and it just sort of jostled me a little bit because we''re so similar'
$ws.Range("F4").Value = 'she''s overweight?'

$ws.Range("D5").Value = 'very fun all right and how are things going in school for you'
$ws.Range("E5").Value = 'pretty good grades have been pretty all right still like 3.7 ish'
$ws.Range("F5").Value = 'well thanks for taking time out to talk to us today'

$ws.Range("D6").Value = 'mm-hmm and do you feel like that''s working pretty well for you or you sometimes exceed that limit or or what do you think'
$ws.Range("E6").Value = 'for the most part it works pretty well because like I said people you know people come asking me trying to pressure me I''ll be like no I''m good I have you know I still have so I''m like I''ll get some later and then never end up getting anything but sometimes they''ll be like hey let''s take a shot and everyone''s got it around like the kitchen or something ready to take a shot and then at that point I''ll be like oh well maybe I can take one more and that''s when I start you know falling into'
$ws.Range("F6").Value = 'so you just want to make sure you''re taking enough?'

$ws.Range("D7").Value = 'uh-huh well why did you say six or seven rather than like a four'
$ws.Range("E7").Value = 'um because I still think it''s it''s important'
$ws.Range("F7").Value = 'Okay so you know you can''t smoke on the weekends, right?'

$ws.Range("D8").Value = 'okay so what''s been going on what is gonna happen'
$ws.Range("E8").Value = 'um well I am gonna have to go to court soon and I know that I''m probably gonna have to get some type of like treatment or something like that but I really just don''t want to so I just I''m just coming here because I don''t know maybe I can like get out of it if I''m getting some type of treatment but I don''t know'
$ws.Range("F8").Value = 'ok so how long will it take?'

$ws.Range("D9").Value = 'well if you don''t mind John if we can go over to the side area and I can discuss with you further kind of see your ankle John'
$ws.Range("E9").Value = 'sure'
$ws.Range("F9").Value = 'okay cool thanks man'

$ws.Range("D10").Value = 'yeah so your financial stress may be increased by the fact that you''re drinking'
$ws.Range("E10").Value = 'yeah kind of defeating the purpose there'
$ws.Range("F10").Value = 'but you''ve got to do something about it right?'

$ws.Range("D11").Value = 'huh okay so on the one side some of the fun things associated with alcohol are people are just more relaxed sort of carefree and kind of let loose a little bit more on the other hand there''s a little bit more I guess planning involved because you have to make sure that you''re doing it safely or as safe as you can and then there can be some other consequences like throwing up or even if it''s not you that sometimes you''re stuck your your fun is sort of ruined when you''re stuck dealing with a friend you drink too much okay any other not-so-good things about the alcohol'
$ws.Range("E11").Value = 'yeah people can get a little like sloppy I guess like they''re all loose like with people and corners doing things that should be done in private and that kind of stuff'
$ws.Range("F11").Value = 'oh wow that''s interesting'

$ws.Range("D12").Value = 'like forget about problems'
$ws.Range("E12").Value = 'yeah'
$ws.Range("F12").Value = 'well let''s talk about this problem'

$ws.Range("D13").Value = 'I want you have a seat here sure so to start off I''m just going to verify your information real quick so we said your name''s Larry seedorf is your date of birth let''s see'
$ws.Range("E13").Value = '129 68 okay'
$ws.Range("F13").Value = 'oh wow that''s pretty cool'

$ws.Range("D14").Value = 'yeah I understand it does sound like you''re worried about getting holes in your teeth though there are a few options that you can do to help improve how you consume sugar would you like to hear about those'
$ws.Range("E14").Value = 'well I don''t like getting feelings so that would be great'
$ws.Range("F14").Value = 'oh ok cool so if you were going to try one thing, what would it be?'

$ws.Range("D15").Value = 'oh nice'
$ws.Range("E15").Value = 'I''m a big boy'
$ws.Range("F15").Value = 'so you''re going to start taking care of yourself better?'

$ws.Range("D16").Value = 'agree with my secret'
$ws.Range("E16").Value = 'honestly I don''t know I mean I was swear I don''t have a problem but how am I gonna you know my best friend will talk to me I mean she drinks more than anyone I know and she said I''m going overboard you know she tried to tell me paste want for everyone she has you know so and I don''t want to tell you this but I got another DUI'
$ws.Range("F16").Value = 'oh yeah right'

$ws.Range("D17").Value = 'right you''ve got it okay good what questions do you have for me about this medication'
$ws.Range("E17").Value = 'I mean you know I mean I just hope that doesn''t cause the pain I mean obviously you know it''s frustrating I''m good not have any problems and then start taking medicine and starting to have problems'
$ws.Range("F17").Value = 'oh ok cool cool cool cool cool cool cool cool cool cool cool cool cool cool cool cool cool cool cool cool cool cool cool cool cool cool cool cool cool cool cool cool cool cool cool cool cool cool cool cool cool cool cool cool cool cool cool cool cool cool cool cool cool cool cool cool cool cool cool cool cool cool cool cool cool cool cool cool cool cool cool cool cool cool cool cool cool cool cool cool cool cool cool cool cool cool cool cool cool cool cool cool cool cool cool cool cool cool cool cool cool cool cool cool cool cool cool cool cool cool cool cool cool cool cool cool cool cool cool cool cool cool cool cool cool cool cool cool cool cool cool cool cool cool cool cool cool cool cool cool cool cool cool cool cool cool cool cool cool cool cool cool cool cool cool cool cool cool cool cool cool cool cool cool cool cool cool cool cool cool cool cool cool cool cool cool cool cool cool cool cool cool cool cool cool cool cool cool cool cool cool cool cool'

$ws.Range("D18").Value = 'and when you use the opposites of a happening and what I want you to do is first put what preceded what that was happening before you decided to use it for this technique so whether it was an argument your mother or somebody cut you off you know you''re getting stressed about a homework assignment whatever happens to be what happened before which behavior you choose to use which skill set you''re going to implement and then how do you sell afterwards okay and it made you feel calmer happier got along better with your mother you were able to finish your assignment without anxiety that whatever that happened to be'
$ws.Range("E18").Value = 'okay'
$ws.Range("F18").Value = 'so you''ve done all these things in order to make yourself feel calmer, more relaxed, less anxious, and more comfortable?'

$ws.Range("D19").Value = 'it''s good to see you again'
$ws.Range("E19").Value = 'it''s nice to see you'
$ws.Range("F19").Value = 'thank you for coming back to visit us'

$ws.Range("D20").Value = 'okay tell me a little bit about how smoking fits in your day it''ll help with the'
$ws.Range("E20").Value = 'I mean I suppose a little less than a pack a day it''s not too bad um the I probably smoked for the last five years but I I mean I could quit if I wanted to but I don''t really want to right now'
$ws.Range("F20").Value = 'oh ok cool'

$ws.Range("D21").Value = 'okay what made you say two instead of one'
$ws.Range("E21").Value = 'well I mean anything''s possible I mean I I guess I could maybe be convinced yeah is there something I mean I are there things I can try that would make me like want to quit smoking'
$ws.Range("F21").Value = 'no, just stop smoking. You know, we all need help in our lives. We should never give up trying to live life as best as we can. It doesn''t matter whether or not you succeed at doing this. If you fail then you failed. That''s it. So quit smoking now. Quit smoking now. Stop smoking now.'
